$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new row 93 (shift old row 93 "note" row down to row 94) ---
# This also naturally clones the number-format styles (s="7"/"8") of the row
# above into the new row 93 cells, matching the source workbook's pattern.
$ws.Rows.Item(93).Insert(-4121)

# --- Fill in the newly-available daily figures for 2020-04-26 (row 92) ---
$ws.Range("B92").Value = 514
$ws.Range("C92").Value = 30028
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 6458

# --- New date row for 2020-04-27 (row 93); B93:E93 stay blank for now ---
$ws.Range("A93").Value = 43948

# --- View: zoom to 85% (pageBreakPreview / sheet layout view) ---
$win = $excel.ActiveWindow
$win.Zoom = 85

# --- Selection moves from B92 to A92 ---
$ws.Range("A92").Select()

# --- Update the workbook-level Print_Area named range ---
$nm = $wb.Names.Item(1)
$nm.RefersTo = "=" + $ws.Name() + "!`$A`$1:`$E`$96"
